# Update the "Burn Down Chart Table" worksheet with the latest hours logged.
# (Formulas in column AC and in the totals row 30 recalculate automatically,
# and the embedded chart on the "Burn Down Chart" sheet reads its series from
# row 30, so it reflects the new numbers as well.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burn Down Chart Table")

$ws.Range("D2").Value  = 2.5
$ws.Range("D3").Value  = 0.5
$ws.Range("D5").Value  = 1
$ws.Range("D6").Value  = 1
$ws.Range("D7").Value  = 1
$ws.Range("F8").Value  = 2
$ws.Range("E9").Value  = 1
$ws.Range("G10").Value = 1
$ws.Range("G11").Value = 0.5
$ws.Range("K19").Value = 1
$ws.Range("M19").Value = 1
$ws.Range("K23").Value = 1

# Leave the selection the way it was left when the file was last saved.
$ws.Range("C22:C23").Select() | Out-Null
